# Apply scheduled-runner value updates to the Bahamut_Profits workbook.
# Each block below rewrites the H:N "profit" columns for one leve row on one
# Job sheet (ALC/ARM/BSM/CRP/CUL/LTW/WVR), reflecting refreshed market prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: update H15, I15, K15, M15
$ws.Cells.Item(15, 8).Value = 94.03
$ws.Cells.Item(15, 9).Value = 94.03
$ws.Cells.Item(15, 11).Value = 282.09
$ws.Cells.Item(15, 13).Value = -113.09

# Row 43: update H43, I43, J43, K43, L43, M43, N43
$ws.Cells.Item(43, 8).Value = 2677
$ws.Cells.Item(43, 9).Value = 4292.5
$ws.Cells.Item(43, 10).Value = 1600
$ws.Cells.Item(43, 11).Value = 4292.5
$ws.Cells.Item(43, 12).Value = 1600
$ws.Cells.Item(43, 13).Value = -4223.5
$ws.Cells.Item(43, 14).Value = -1738

# Row 53: update H53, J53, L53, N53
$ws.Cells.Item(53, 8).Value = 278.25
$ws.Cells.Item(53, 10).Value = 649.6
$ws.Cells.Item(53, 12).Value = 649.6
$ws.Cells.Item(53, 14).Value = -1923.6

# Row 62: update H62, I62, J62, K62, L62, M62, N62
$ws.Cells.Item(62, 8).Value = 53168.24
$ws.Cells.Item(62, 9).Value = 71208.92999999999
$ws.Cells.Item(62, 10).Value = 8066.5
$ws.Cells.Item(62, 11).Value = 71208.92999999999
$ws.Cells.Item(62, 12).Value = 8066.5
$ws.Cells.Item(62, 13).Value = -70584.92999999999
$ws.Cells.Item(62, 14).Value = -9314.5

# Row 65: update H65, I65, J65, K65, L65, M65, N65
$ws.Cells.Item(65, 8).Value = 53168.24
$ws.Cells.Item(65, 9).Value = 71208.92999999999
$ws.Cells.Item(65, 10).Value = 8066.5
$ws.Cells.Item(65, 11).Value = 356044.65
$ws.Cells.Item(65, 12).Value = 40332.5
$ws.Cells.Item(65, 13).Value = -352924.65
$ws.Cells.Item(65, 14).Value = -46572.5

# Row 103: update H103, I103, J103, K103, L103, M103, N103
$ws.Cells.Item(103, 8).Value = 704.6
$ws.Cells.Item(103, 9).Value = 802.5333000000001
$ws.Cells.Item(103, 10).Value = 410.8
$ws.Cells.Item(103, 11).Value = 2407.5999
$ws.Cells.Item(103, 12).Value = 1232.4
$ws.Cells.Item(103, 13).Value = -1821.5999
$ws.Cells.Item(103, 14).Value = -2404.4

# Row 116: update H116, I116, J116, K116, L116, M116, N116
$ws.Cells.Item(116, 8).Value = 3742.8928
$ws.Cells.Item(116, 9).Value = 3695.739
$ws.Cells.Item(116, 10).Value = 3959.8
$ws.Cells.Item(116, 11).Value = 3695.739
$ws.Cells.Item(116, 12).Value = 3959.8
$ws.Cells.Item(116, 13).Value = -253.739
$ws.Cells.Item(116, 14).Value = -10843.8

# Row 125: update H125, I125, J125, K125, L125, M125, N125
$ws.Cells.Item(125, 8).Value = 1550
$ws.Cells.Item(125, 9).Value = 1500
$ws.Cells.Item(125, 10).Value = 1600
$ws.Cells.Item(125, 11).Value = 13500
$ws.Cells.Item(125, 12).Value = 14400
$ws.Cells.Item(125, 13).Value = -11040
$ws.Cells.Item(125, 14).Value = -19320

# Row 129: update H129, J129, L129, N129
$ws.Cells.Item(129, 8).Value = 882895.0600000001
$ws.Cells.Item(129, 10).Value = 1029968.75
$ws.Cells.Item(129, 12).Value = 3089906.25
$ws.Cells.Item(129, 14).Value = -3099906.25

# Row 132: update H132, I132, J132, K132, L132, M132, N132
$ws.Cells.Item(132, 8).Value = 1888676.6
$ws.Cells.Item(132, 9).Value = 1661.579
$ws.Cells.Item(132, 10).Value = 6669114.5
$ws.Cells.Item(132, 11).Value = 4984.737
$ws.Cells.Item(132, 12).Value = 20007343.5
$ws.Cells.Item(132, 13).Value = -2454.737
$ws.Cells.Item(132, 14).Value = -20012403.5

$ws = $wb.Worksheets.Item("ARM")
# Row 61: update H61, I61, J61, K61, L61, M61, N61
$ws.Cells.Item(61, 8).Value = 2753.64
$ws.Cells.Item(61, 9).Value = 3015.8667
$ws.Cells.Item(61, 10).Value = 2360.3
$ws.Cells.Item(61, 11).Value = 3015.8667
$ws.Cells.Item(61, 12).Value = 2360.3
$ws.Cells.Item(61, 13).Value = -2803.8667
$ws.Cells.Item(61, 14).Value = -2784.3

# Row 74: update H74, I74, J74, K74, L74, M74, N74
$ws.Cells.Item(74, 8).Value = 1026.3405
$ws.Cells.Item(74, 9).Value = 997.6111
$ws.Cells.Item(74, 10).Value = 1120.3636
$ws.Cells.Item(74, 11).Value = 997.6111
$ws.Cells.Item(74, 12).Value = 1120.3636
$ws.Cells.Item(74, 13).Value = -123.6111
$ws.Cells.Item(74, 14).Value = -2868.3636

# Row 77: update H77, I77, J77, K77, L77, M77, N77
$ws.Cells.Item(77, 8).Value = 1026.3405
$ws.Cells.Item(77, 9).Value = 997.6111
$ws.Cells.Item(77, 10).Value = 1120.3636
$ws.Cells.Item(77, 11).Value = 4988.055499999999
$ws.Cells.Item(77, 12).Value = 5601.817999999999
$ws.Cells.Item(77, 13).Value = -620.0554999999995
$ws.Cells.Item(77, 14).Value = -14337.818

# Row 102: update H102, I102, J102, K102, L102, M102, N102
$ws.Cells.Item(102, 8).Value = 5300
$ws.Cells.Item(102, 9).Value = 5560
$ws.Cells.Item(102, 10).Value = 4000
$ws.Cells.Item(102, 11).Value = 5560
$ws.Cells.Item(102, 12).Value = 4000
$ws.Cells.Item(102, 13).Value = -3938
$ws.Cells.Item(102, 14).Value = -7244

# Row 103: update H103, J103, L103, N103
$ws.Cells.Item(103, 8).Value = 40000
$ws.Cells.Item(103, 10).Value = 40000
$ws.Cells.Item(103, 12).Value = 40000
$ws.Cells.Item(103, 14).Value = -42344

# Row 136: update H136, I136, J136, K136, L136, M136, N136
$ws.Cells.Item(136, 8).Value = 2753.64
$ws.Cells.Item(136, 9).Value = 3015.8667
$ws.Cells.Item(136, 10).Value = 2360.3
$ws.Cells.Item(136, 11).Value = 9047.6001
$ws.Cells.Item(136, 12).Value = 7080.900000000001
$ws.Cells.Item(136, 13).Value = -6497.6001
$ws.Cells.Item(136, 14).Value = -12180.9

$ws = $wb.Worksheets.Item("BSM")
# Row 97: update H97, I97, K97, M97
$ws.Cells.Item(97, 8).Value = 2322.2856
$ws.Cells.Item(97, 9).Value = 2322.2856
$ws.Cells.Item(97, 11).Value = 2322.2856
$ws.Cells.Item(97, 13).Value = -1331.2856

# Row 103: update H103, J103, L103, N103
$ws.Cells.Item(103, 8).Value = 25550
$ws.Cells.Item(103, 10).Value = 25550
$ws.Cells.Item(103, 12).Value = 25550
$ws.Cells.Item(103, 14).Value = -27894

$ws = $wb.Worksheets.Item("CRP")
# Row 31: update H31, I31, J31, K31, L31, M31, N31
$ws.Cells.Item(31, 8).Value = 33213.883
$ws.Cells.Item(31, 9).Value = 3953.762
$ws.Cells.Item(31, 10).Value = 80480.234
$ws.Cells.Item(31, 11).Value = 3953.762
$ws.Cells.Item(31, 12).Value = 80480.234
$ws.Cells.Item(31, 13).Value = -3658.762
$ws.Cells.Item(31, 14).Value = -81070.234

# Row 34: update H34, I34, J34, K34, L34, M34, N34
$ws.Cells.Item(34, 8).Value = 33213.883
$ws.Cells.Item(34, 9).Value = 3953.762
$ws.Cells.Item(34, 10).Value = 80480.234
$ws.Cells.Item(34, 11).Value = 3953.762
$ws.Cells.Item(34, 12).Value = 80480.234
$ws.Cells.Item(34, 13).Value = -3751.762
$ws.Cells.Item(34, 14).Value = -80884.234

# Row 97: update H97, J97, L97, N97
$ws.Cells.Item(97, 8).Value = 19700
$ws.Cells.Item(97, 10).Value = 19700
$ws.Cells.Item(97, 12).Value = 19700
$ws.Cells.Item(97, 14).Value = -21682

# Row 99: update H99, I99, J99, K99, L99, M99, N99
$ws.Cells.Item(99, 8).Value = 2993.5293
$ws.Cells.Item(99, 9).Value = 2936.2964
$ws.Cells.Item(99, 10).Value = 3214.2856
$ws.Cells.Item(99, 11).Value = 2936.2964
$ws.Cells.Item(99, 12).Value = 3214.2856
$ws.Cells.Item(99, 13).Value = -1438.2964
$ws.Cells.Item(99, 14).Value = -6210.2856

# Row 126: update H126, I126, J126, K126, L126, M126, N126
$ws.Cells.Item(126, 8).Value = 2993.5293
$ws.Cells.Item(126, 9).Value = 2936.2964
$ws.Cells.Item(126, 10).Value = 3214.2856
$ws.Cells.Item(126, 11).Value = 8808.889200000001
$ws.Cells.Item(126, 12).Value = 9642.856800000001
$ws.Cells.Item(126, 13).Value = -6338.889200000001
$ws.Cells.Item(126, 14).Value = -14582.8568

$ws = $wb.Worksheets.Item("CUL")
# Row 5: update H5, I5, J5, K5, L5, M5, N5
$ws.Cells.Item(5, 8).Value = 2177.2766
$ws.Cells.Item(5, 9).Value = 1250.375
$ws.Cells.Item(5, 10).Value = 2655.6775
$ws.Cells.Item(5, 11).Value = 3751.125
$ws.Cells.Item(5, 12).Value = 7967.032499999999
$ws.Cells.Item(5, 13).Value = -3639.125
$ws.Cells.Item(5, 14).Value = -8191.032499999999

# Row 121: update H121, I121, J121, K121, L121, M121, N121
$ws.Cells.Item(121, 8).Value = 828.65
$ws.Cells.Item(121, 9).Value = 350
$ws.Cells.Item(121, 10).Value = 853.8421
$ws.Cells.Item(121, 11).Value = 1050
$ws.Cells.Item(121, 12).Value = 2561.5263
$ws.Cells.Item(121, 13).Value = 260
$ws.Cells.Item(121, 14).Value = -5181.5263

# Row 132: update H132, I132, J132, K132, L132, M132, N132
$ws.Cells.Item(132, 8).Value = 1053.7727
$ws.Cells.Item(132, 9).Value = 939.5625
$ws.Cells.Item(132, 10).Value = 1358.3334
$ws.Cells.Item(132, 11).Value = 8456.0625
$ws.Cells.Item(132, 12).Value = 12225.0006
$ws.Cells.Item(132, 13).Value = -5926.0625
$ws.Cells.Item(132, 14).Value = -17285.0006

# Row 135: update H135, I135, J135, K135, L135, M135, N135
$ws.Cells.Item(135, 8).Value = 2177.2766
$ws.Cells.Item(135, 9).Value = 1250.375
$ws.Cells.Item(135, 10).Value = 2655.6775
$ws.Cells.Item(135, 11).Value = 11253.375
$ws.Cells.Item(135, 12).Value = 23901.0975
$ws.Cells.Item(135, 13).Value = -8718.375
$ws.Cells.Item(135, 14).Value = -28971.0975

# Row 137: update H137, J137, L137, N137
$ws.Cells.Item(137, 8).Value = 51203.57
$ws.Cells.Item(137, 10).Value = 115388.336
$ws.Cells.Item(137, 12).Value = 346165.008
$ws.Cells.Item(137, 14).Value = -356365.008

$ws = $wb.Worksheets.Item("LTW")
# Row 7: update H7, I7, J7, K7, L7, M7, N7
$ws.Cells.Item(7, 8).Value = 2651.64
$ws.Cells.Item(7, 9).Value = 2574
$ws.Cells.Item(7, 10).Value = 2851.2856
$ws.Cells.Item(7, 11).Value = 2574
$ws.Cells.Item(7, 12).Value = 2851.2856
$ws.Cells.Item(7, 13).Value = -2462
$ws.Cells.Item(7, 14).Value = -3075.2856

# Row 68: update H68, I68, J68, K68, L68, M68, N68
$ws.Cells.Item(68, 8).Value = 2652.7036
$ws.Cells.Item(68, 9).Value = 2540.1667
$ws.Cells.Item(68, 10).Value = 2877.7778
$ws.Cells.Item(68, 11).Value = 2540.1667
$ws.Cells.Item(68, 12).Value = 2877.7778
$ws.Cells.Item(68, 13).Value = -1791.1667
$ws.Cells.Item(68, 14).Value = -4375.7778

# Row 71: update H71, I71, J71, K71, L71, M71, N71
$ws.Cells.Item(71, 8).Value = 2652.7036
$ws.Cells.Item(71, 9).Value = 2540.1667
$ws.Cells.Item(71, 10).Value = 2877.7778
$ws.Cells.Item(71, 11).Value = 12700.8335
$ws.Cells.Item(71, 12).Value = 14388.889
$ws.Cells.Item(71, 13).Value = -8956.833500000001
$ws.Cells.Item(71, 14).Value = -21876.889

# Row 126: update H126, I126, J126, K126, L126, M126, N126
$ws.Cells.Item(126, 8).Value = 2651.64
$ws.Cells.Item(126, 9).Value = 2574
$ws.Cells.Item(126, 10).Value = 2851.2856
$ws.Cells.Item(126, 11).Value = 7722
$ws.Cells.Item(126, 12).Value = 8553.856800000001
$ws.Cells.Item(126, 13).Value = -5252
$ws.Cells.Item(126, 14).Value = -13493.8568

$ws = $wb.Worksheets.Item("WVR")
# Row 93: update H93, J93, L93, N93
$ws.Cells.Item(93, 8).Value = 20896.572
$ws.Cells.Item(93, 10).Value = 20896.572
$ws.Cells.Item(93, 12).Value = 20896.572
$ws.Cells.Item(93, 14).Value = -25888.572

# Row 122: update H122, I122, J122, K122, L122, M122, N122
$ws.Cells.Item(122, 8).Value = 90910024
$ws.Cells.Item(122, 9).Value = 142858000
$ws.Cells.Item(122, 10).Value = 1075
$ws.Cells.Item(122, 11).Value = 428574000
$ws.Cells.Item(122, 12).Value = 3225
$ws.Cells.Item(122, 13).Value = -428571550
$ws.Cells.Item(122, 14).Value = -8125
